$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.369957804679871
$ws.Range("B1").Value = 2.838464498519897
$ws.Range("C1").Value = 3.793054342269897
$ws.Range("D1").Value = 3.646950006484985
$ws.Range("E1").Value = 0.9798374176025391
